$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = [double]"1.168788"
$ws.Range("H2").Value2 = [double]"3.506364"
$ws.Range("I2").Value2 = [double]"0.0139255825297802"
$ws.Range("J2").Value2 = [double]"0.01408364098536324"
$ws.Range("M2").Value2 = [double]"0.169654"
$ws.Range("N2").Value2 = [double]"0.508962"
$ws.Range("O2").Value2 = [double]"0.006094264463659866"
$ws.Range("P2").Value2 = [double]"0.006534681579452628"
$ws.Range("Q2").Value2 = [double]"0.198289559352"
$ws.Range("R2").Value2 = [double]"1.784606034168"
$ws.Range("S2").Value2 = [double]"8.486618274700214E-05"
$ws.Range("T2").Value2 = [double]"9.203210931867721E-05"
$ws.Range("G3").Value2 = [double]"1.168788"
$ws.Range("H3").Value2 = [double]"3.506364"
$ws.Range("I3").Value2 = [double]"0.0139255825297802"
$ws.Range("J3").Value2 = [double]"0.01408364098536324"
$ws.Range("O3").Value2 = [double]"0.7900017288527916"
$ws.Range("P3").Value2 = [double]"0.8470931604713817"
$ws.Range("Q3").Value2 = [double]"25.70434802028"
$ws.Range("R3").Value2 = [double]"231.33913218252"
$ws.Range("S3").Value2 = [double]"0.01100123427380859"
$ws.Range("T3").Value2 = [double]"0.01193015595323563"
$ws.Range("G4").Value2 = [double]"1.168788"
$ws.Range("H4").Value2 = [double]"3.506364"
$ws.Range("I4").Value2 = [double]"0.0139255825297802"
$ws.Range("J4").Value2 = [double]"0.01408364098536324"
$ws.Range("M4").Value2 = [double]"0.04769766666666667"
$ws.Range("N4").Value2 = [double]"0.143093"
$ws.Range("O4").Value2 = [double]"0.001713382501834088"
$ws.Range("P4").Value2 = [double]"0.001837204332049573"
$ws.Range("Q4").Value2 = [double]"0.055748460428"
$ws.Range("R4").Value2 = [double]"0.501736143852"
$ws.Range("S4").Value2 = [double]"2.385984943437187E-05"
$ws.Range("T4").Value2 = [double]"2.587452622934026E-05"
$ws.Range("G5").Value2 = [double]"1.168788"
$ws.Range("H5").Value2 = [double]"3.506364"
$ws.Range("I5").Value2 = [double]"0.0139255825297802"
$ws.Range("J5").Value2 = [double]"0.01408364098536324"
$ws.Range("M5").Value2 = [double]"5.6286445"
$ws.Range("N5").Value2 = [double]"11.257289"
$ws.Range("O5").Value2 = [double]"0.2021906241817143"
$ws.Range("P5").Value2 = [double]"0.1445349536171162"
$ws.Range("Q5").Value2 = [double]"6.578692147866"
$ws.Range("R5").Value2 = [double]"39.472152887196"
$ws.Range("S5").Value2 = [double]"0.002815622223790235"
$ws.Range("T5").Value2 = [double]"0.002035578396579592"
$ws.Range("I6").Value2 = [double]"0.9504675141158767"
$ws.Range("J6").Value2 = [double]"0.9612555315680539"
$ws.Range("M6").Value2 = [double]"0.169654"
$ws.Range("N6").Value2 = [double]"0.508962"
$ws.Range("O6").Value2 = [double]"0.006094264463659866"
$ws.Range("P6").Value2 = [double]"0.006534681579452628"
$ws.Range("Q6").Value2 = [double]"13.53392464188733"
$ws.Range("R6").Value2 = [double]"121.805321776986"
$ws.Range("S6").Value2 = [double]"0.00579240039513952"
$ws.Range("T6").Value2 = [double]"0.006281498815284706"
$ws.Range("I7").Value2 = [double]"0.9504675141158767"
$ws.Range("J7").Value2 = [double]"0.9612555315680539"
$ws.Range("O7").Value2 = [double]"0.7900017288527916"
$ws.Range("P7").Value2 = [double]"0.8470931604713817"
$ws.Range("S7").Value2 = [double]"0.7508709793699577"
$ws.Range("T7").Value2 = [double]"0.8142729862565808"
$ws.Range("I8").Value2 = [double]"0.9504675141158767"
$ws.Range("J8").Value2 = [double]"0.9612555315680539"
$ws.Range("M8").Value2 = [double]"0.04769766666666667"
$ws.Range("N8").Value2 = [double]"0.143093"
$ws.Range("O8").Value2 = [double]"0.001713382501834088"
$ws.Range("P8").Value2 = [double]"0.001837204332049573"
$ws.Range("Q8").Value2 = [double]"3.805018604103222"
$ws.Range("R8").Value2 = [double]"34.245167436929"
$ws.Range("S8").Value2 = [double]"0.001628514407247887"
$ws.Range("T8").Value2 = [double]"0.001766022826803444"
$ws.Range("I9").Value2 = [double]"0.9504675141158767"
$ws.Range("J9").Value2 = [double]"0.9612555315680539"
$ws.Range("M9").Value2 = [double]"5.6286445"
$ws.Range("N9").Value2 = [double]"11.257289"
$ws.Range("O9").Value2 = [double]"0.2021906241817143"
$ws.Range("P9").Value2 = [double]"0.1445349536171162"
$ws.Range("Q9").Value2 = [double]"449.0177095675529"
$ws.Range("R9").Value2 = [double]"2694.106257405317"
$ws.Range("S9").Value2 = [double]"0.1921756199435315"
$ws.Range("T9").Value2 = [double]"0.138935023669385"
$ws.Range("E10").Value2 = [double]"3"
$ws.Range("F10").Value2 = [double]"1"
$ws.Range("G10").Value2 = [double]"0.09388066666666665"
$ws.Range("H10").Value2 = [double]"0.2816419999999999"
$ws.Range("I10").Value2 = [double]"0.001118545854010694"
$ws.Range("J10").Value2 = [double]"0.001131241597962925"
$ws.Range("M10").Value2 = [double]"0.169654"
$ws.Range("N10").Value2 = [double]"0.508962"
$ws.Range("O10").Value2 = [double]"0.006094264463659866"
$ws.Range("P10").Value2 = [double]"0.006534681579452628"
$ws.Range("Q10").Value2 = [double]"0.01592723062266667"
$ws.Range("R10").Value2 = [double]"0.143345075604"
$ws.Range("S10").Value2 = [double]"6.816714249071452E-06"
$ws.Range("T10").Value2 = [double]"7.392303632118879E-06"
$ws.Range("E11").Value2 = [double]"3"
$ws.Range("F11").Value2 = [double]"1"
$ws.Range("G11").Value2 = [double]"0.09388066666666665"
$ws.Range("H11").Value2 = [double]"0.2816419999999999"
$ws.Range("I11").Value2 = [double]"0.001118545854010694"
$ws.Range("J11").Value2 = [double]"0.001131241597962925"
$ws.Range("O11").Value2 = [double]"0.7900017288527916"
$ws.Range("P11").Value2 = [double]"0.8470931604713817"
$ws.Range("Q11").Value2 = [double]"2.06465272434"
$ws.Range("R11").Value2 = [double]"18.58187451905999"
$ws.Range("S11").Value2 = [double]"0.0008836531584695709"
$ws.Range("T11").Value2 = [double]"0.0009582670204751101"
$ws.Range("E12").Value2 = [double]"3"
$ws.Range("F12").Value2 = [double]"1"
$ws.Range("G12").Value2 = [double]"0.09388066666666665"
$ws.Range("H12").Value2 = [double]"0.2816419999999999"
$ws.Range("I12").Value2 = [double]"0.001118545854010694"
$ws.Range("J12").Value2 = [double]"0.001131241597962925"
$ws.Range("M12").Value2 = [double]"0.04769766666666667"
$ws.Range("N12").Value2 = [double]"0.143093"
$ws.Range("O12").Value2 = [double]"0.001713382501834088"
$ws.Range("P12").Value2 = [double]"0.001837204332049573"
$ws.Range("Q12").Value2 = [double]"0.00447788874511111"
$ws.Range("R12").Value2 = [double]"0.04030099870599999"
$ws.Range("S12").Value2 = [double]"1.916496893760991E-06"
$ws.Range("T12").Value2 = [double]"2.078321964372166E-06"
$ws.Range("E13").Value2 = [double]"3"
$ws.Range("F13").Value2 = [double]"1"
$ws.Range("G13").Value2 = [double]"0.09388066666666665"
$ws.Range("H13").Value2 = [double]"0.2816419999999999"
$ws.Range("I13").Value2 = [double]"0.001118545854010694"
$ws.Range("J13").Value2 = [double]"0.001131241597962925"
$ws.Range("M13").Value2 = [double]"5.6286445"
$ws.Range("N13").Value2 = [double]"11.257289"
$ws.Range("O13").Value2 = [double]"0.2021906241817143"
$ws.Range("P13").Value2 = [double]"0.1445349536171162"
$ws.Range("Q13").Value2 = [double]"0.5284208980896666"
$ws.Range("R13").Value2 = [double]"3.170525388538"
$ws.Range("S13").Value2 = [double]"0.000226159484398291"
$ws.Range("T13").Value2 = [double]"0.0001635039518913237"
$ws.Range("G14").Value2 = [double]"2.8258325"
$ws.Range("H14").Value2 = [double]"5.651664999999999"
$ws.Range("I14").Value2 = [double]"0.03366852131788238"
$ws.Range("J14").Value2 = [double]"0.0227004443433548"
$ws.Range("M14").Value2 = [double]"0.169654"
$ws.Range("N14").Value2 = [double]"0.508962"
$ws.Range("O14").Value2 = [double]"0.006094264463659866"
$ws.Range("P14").Value2 = [double]"0.006534681579452628"
$ws.Range("Q14").Value2 = [double]"0.479413786955"
$ws.Range("R14").Value2 = [double]"2.87648272173"
$ws.Range("S14").Value2 = [double]"0.0002051848730115452"
$ws.Range("T14").Value2 = [double]"0.0001483401754959102"
$ws.Range("G15").Value2 = [double]"2.8258325"
$ws.Range("H15").Value2 = [double]"5.651664999999999"
$ws.Range("I15").Value2 = [double]"0.03366852131788238"
$ws.Range("J15").Value2 = [double]"0.0227004443433548"
$ws.Range("O15").Value2 = [double]"0.7900017288527916"
$ws.Range("P15").Value2 = [double]"0.8470931604713817"
$ws.Range("Q15").Value2 = [double]"62.146584348075"
$ws.Range("R15").Value2 = [double]"372.8795060884499"
$ws.Range("S15").Value2 = [double]"0.02659819004904415"
$ws.Range("T15").Value2 = [double]"0.01922939114291712"
$ws.Range("G16").Value2 = [double]"2.8258325"
$ws.Range("H16").Value2 = [double]"5.651664999999999"
$ws.Range("I16").Value2 = [double]"0.03366852131788238"
$ws.Range("J16").Value2 = [double]"0.0227004443433548"
$ws.Range("M16").Value2 = [double]"0.04769766666666667"
$ws.Range("N16").Value2 = [double]"0.143093"
$ws.Range("O16").Value2 = [double]"0.001713382501834088"
$ws.Range("P16").Value2 = [double]"0.001837204332049573"
$ws.Range("Q16").Value2 = [double]"0.1347856166408333"
$ws.Range("R16").Value2 = [double]"0.808713699845"
$ws.Range("S16").Value2 = [double]"5.768705528868765E-05"
$ws.Range("T16").Value2 = [double]"4.170535468706167E-05"
$ws.Range("G17").Value2 = [double]"2.8258325"
$ws.Range("H17").Value2 = [double]"5.651664999999999"
$ws.Range("I17").Value2 = [double]"0.03366852131788238"
$ws.Range("J17").Value2 = [double]"0.0227004443433548"
$ws.Range("M17").Value2 = [double]"5.6286445"
$ws.Range("N17").Value2 = [double]"11.257289"
$ws.Range("O17").Value2 = [double]"0.2021906241817143"
$ws.Range("P17").Value2 = [double]"0.1445349536171162"
$ws.Range("Q17").Value2 = [double]"15.90560655904625"
$ws.Range("R17").Value2 = [double]"63.62242623618499"
$ws.Range("S17").Value2 = [double]"0.006807459340537992"
$ws.Range("T17").Value2 = [double]"0.003281007670254714"
$ws.Range("E18").Value2 = [double]"2"
$ws.Range("F18").Value2 = [double]"0.6666666666666666"
$ws.Range("G18").Value2 = [double]"0.06880966666666666"
$ws.Range("H18").Value2 = [double]"0.206429"
$ws.Range("I18").Value2 = [double]"0.0008198361824499672"
$ws.Range("J18").Value2 = [double]"0.0008291415052651543"
$ws.Range("M18").Value2 = [double]"0.169654"
$ws.Range("N18").Value2 = [double]"0.508962"
$ws.Range("O18").Value2 = [double]"0.006094264463659866"
$ws.Range("P18").Value2 = [double]"0.006534681579452628"
$ws.Range("Q18").Value2 = [double]"0.01167383518866666"
$ws.Range("R18").Value2 = [double]"0.105064516698"
$ws.Range("S18").Value2 = [double]"4.996298512727402E-06"
$ws.Range("T18").Value2 = [double]"5.418175721215828E-06"
$ws.Range("E19").Value2 = [double]"2"
$ws.Range("F19").Value2 = [double]"0.6666666666666666"
$ws.Range("G19").Value2 = [double]"0.06880966666666666"
$ws.Range("H19").Value2 = [double]"0.206429"
$ws.Range("I19").Value2 = [double]"0.0008198361824499672"
$ws.Range("J19").Value2 = [double]"0.0008291415052651543"
$ws.Range("O19").Value2 = [double]"0.7900017288527916"
$ws.Range("P19").Value2 = [double]"0.8470931604713817"
$ws.Range("Q19").Value2 = [double]"1.51328352033"
$ws.Range("R19").Value2 = [double]"13.61955168297"
$ws.Range("S19").Value2 = [double]"0.0006476720015115468"
$ws.Range("T19").Value2 = [double]"0.0007023600981730583"
$ws.Range("E20").Value2 = [double]"2"
$ws.Range("F20").Value2 = [double]"0.6666666666666666"
$ws.Range("G20").Value2 = [double]"0.06880966666666666"
$ws.Range("H20").Value2 = [double]"0.206429"
$ws.Range("I20").Value2 = [double]"0.0008198361824499672"
$ws.Range("J20").Value2 = [double]"0.0008291415052651543"
$ws.Range("M20").Value2 = [double]"0.04769766666666667"
$ws.Range("N20").Value2 = [double]"0.143093"
$ws.Range("O20").Value2 = [double]"0.001713382501834088"
$ws.Range("P20").Value2 = [double]"0.001837204332049573"
$ws.Range("Q20").Value2 = [double]"0.003282060544111111"
$ws.Range("R20").Value2 = [double]"0.029538544897"
$ws.Range("S20").Value2 = [double]"1.404692969380233E-06"
$ws.Range("T20").Value2 = [double]"1.523302365355245E-06"
$ws.Range("E21").Value2 = [double]"2"
$ws.Range("F21").Value2 = [double]"0.6666666666666666"
$ws.Range("G21").Value2 = [double]"0.06880966666666666"
$ws.Range("H21").Value2 = [double]"0.206429"
$ws.Range("I21").Value2 = [double]"0.0008198361824499672"
$ws.Range("J21").Value2 = [double]"0.0008291415052651543"
$ws.Range("M21").Value2 = [double]"5.6286445"
$ws.Range("N21").Value2 = [double]"11.257289"
$ws.Range("O21").Value2 = [double]"0.2021906241817143"
$ws.Range("P21").Value2 = [double]"0.1445349536171162"
$ws.Range("Q21").Value2 = [double]"0.3873051518301666"
$ws.Range("R21").Value2 = [double]"2.323830910981"
$ws.Range("S21").Value2 = [double]"0.0001657631894563127"
$ws.Range("T21").Value2 = [double]"0.000119839929005525"
